# Add a "Cumulate" column (K) with "No" values for the existing data rows,
# matching the investor KPI mapping's new field. Apply an explicit font
# color (kept as automatic/black here) on the new column so a distinct
# cell style gets recorded for it, mirroring the style change introduced
# alongside the "Value Bridge" coloring update in the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "Cumulate"
$ws.Range("K2").Value = "No"
$ws.Range("K3").Value = "No"

$ws.Range("K1:K3").Font.ColorIndex = 1

$ws.Range("K4").Select()
